$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell value updates as described by the diff (crypto price/volume refresh)
# Columns B-E in this sheet are authored as plain text; cells whose new value
# looks like a number are forced to Text format first so Excel does not silently
# convert them (which would lose trailing zeros / thousand-dot formatting).
$ws.Range("D2").Value = "67.192.15"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "3.877.12"
$ws.Range("E3").Value = "  +1.46%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "475.01"
$ws.Range("E5").Value = "  +5.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.31"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.613"
$ws.Range("E7").Value = "  -1.38%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.719"
$ws.Range("E9").Value = "  -2.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("E10").Value = "  +7.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000350"
$ws.Range("E11").Value = "  +10.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.20"
$ws.Range("E12").Value = "  -3.45%  "
$ws.Range("D13").Value = "4.474.54"
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.18"
$ws.Range("E14").Value = "  -1.72%  "
$ws.Range("B15").Value = "Uniswap"
$ws.Range("C15").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.55"
$ws.Range("E15").Value = "  -1.99%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.871.97"
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.67"
$ws.Range("E18").Value = "  -1.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.12"
$ws.Range("E19").Value = "  -4.23%  "
$ws.Range("D20").Value = "67.155.66"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "427.75"
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.31"
$ws.Range("E22").Value = "  +2.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.29"
$ws.Range("E23").Value = "  -3.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.39"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.52"
$ws.Range("E25").Value = "  +1.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "37.82"
$ws.Range("E26").Value = "  +1.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.20"
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.09"
$ws.Range("E28").Value = "  +3.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "724.89"
$ws.Range("E29").Value = "  -1.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.17"
$ws.Range("E30").Value = "  -4.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.128"
$ws.Range("E31").Value = "  -4.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.80"
$ws.Range("E32").Value = "  +2.63%  "
$ws.Range("B33").Value = "PEPE"
$ws.Range("C33").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D33").Value = "0.0₃0886"
$ws.Range("E33").Value = "  +29.76%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "41.78"
$ws.Range("E34").Value = "  -2.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.58"
$ws.Range("E35").Value = "  +2.35%  "
$ws.Range("E36").Value = "  -5.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.42"
$ws.Range("E38").Value = "  -3.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0465"
$ws.Range("E39").Value = "  -2.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.82"
$ws.Range("E40").Value = "  +5.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.02"
$ws.Range("E41").Value = "  +3.72%  "
$ws.Range("E42").Value = "  +11.13%  "
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.343"
$ws.Range("E43").Value = "  -2.68%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.139"
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.41"
$ws.Range("E46").Value = "  -1.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.15"
$ws.Range("E47").Value = "  +1.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "145.73"
$ws.Range("E48").Value = "  +1.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.16"
$ws.Range("E49").Value = "  -3.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.84"
$ws.Range("E50").Value = "  -1.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.57"
$ws.Range("E51").Value = "  -1.80%  "
